$p = $ppt.ActivePresentation

# --- Slide 2: merge the three "DQN (" / "Deep Q-Networks" / ")" runs into one run ---
$s2 = $p.Slides.Item(2)
$dqnShape = $s2.Shapes.Item(4)
$dqnRange = $dqnShape.TextFrame.TextRange
$dqnFound = $dqnRange.Find("DQN (")
$mergeStart = $dqnFound.Start - 1
$mergeRange = $dqnRange.Characters($mergeStart, 22)
$mergeRange.Text = "`tDQN (Deep Q-Networks)"

# --- Slide 3: drop the empty tooltip="" attribute on the two reachable hyperlinks ---
$s3 = $p.Slides.Item(3)
$envShape = $s3.Shapes.Item(2)
$envRange = $envShape.TextFrame.TextRange

$gymFound = $envRange.Find("gym.env")
$gymHyperlink = $gymFound.ActionSettings(1).Hyperlink
$gymHyperlink.ScreenTip = ""

$repoFound = $envRange.Find("https://github.com/YueNing/tn_source_code")
$repoHyperlink = $repoFound.ActionSettings(1).Hyperlink
$repoHyperlink.ScreenTip = ""
